$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.386.27'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.099.10'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.78%  '
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5354'
$ws.Range('E7').Value = '  +3.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4443'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.84'
$ws.Range('E9').Value = '  +3.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09388'
$ws.Range('E10').Value = '  +1.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.172'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.74'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.188.01'
$ws.Range('E13').Value = '  +4.14%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.555'
$ws.Range('E14').Value = '  +3.29%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.909'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '101.63'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001162'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '21.20'
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06680'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.335'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.416.02'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.55'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.315'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.90'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.90'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.521'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.776'
$ws.Range('E29').Value = '  +7.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.81'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.145'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1059'
$ws.Range('E32').Value = '  +0.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.653'
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.268'
$ws.Range('E34').Value = '  +1.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.854'
$ws.Range('E35').Value = '  -2.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.17'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('E37').Value = '  +2.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06813'
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.7024'
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.66'
$ws.Range('E40').Value = '  +1.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.350'
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2223'
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6859'
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.37'
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.343'
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.389'
$ws.Range('E47').Value = '  +19.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.638'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.235'
$ws.Range('E49').Value = '  +9.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000344'
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.220'
$ws.Range('E51').Value = '  +0.04%  '
